# Update the generated two-digit / one-digit division answers in the worksheet table.
# Each data row (1, 5, 9, 13, 17 of the single table) holds 5 answer cells; the new
# values replace the previous ones cell-by-cell, preserving existing formatting by
# writing directly into each cell's Range.Text (Find.Execute is avoided here because
# it searches the whole story rather than the scoped Range in this runtime).
$d = $word.ActiveDocument
$table = $d.Tables(1)

$table.Cell(1, 1).Range.Text = "89÷4=22, 1"  # was "63÷8=7, 7"
$table.Cell(1, 2).Range.Text = "89÷4=22, 1"  # was "28÷4=7, 0"
$table.Cell(1, 3).Range.Text = "56÷9=6, 2"  # was "57÷6=9, 3"
$table.Cell(1, 4).Range.Text = "37÷4=9, 1"  # was "95÷5=19, 0"
$table.Cell(1, 5).Range.Text = "88÷7=12, 4"  # was "36÷7=5, 1"
$table.Cell(5, 1).Range.Text = "19÷2=9, 1"  # was "66÷2=33, 0"
$table.Cell(5, 2).Range.Text = "65÷5=13, 0"  # was "76÷5=15, 1"
$table.Cell(5, 3).Range.Text = "27÷2=13, 1"  # was "93÷5=18, 3"
$table.Cell(5, 4).Range.Text = "23÷9=2, 5"  # was "84÷8=10, 4"
$table.Cell(5, 5).Range.Text = "69÷3=23, 0"  # was "38÷8=4, 6"
$table.Cell(9, 1).Range.Text = "30÷6=5, 0"  # was "34÷5=6, 4"
$table.Cell(9, 2).Range.Text = "83÷4=20, 3"  # was "41÷9=4, 5"
$table.Cell(9, 3).Range.Text = "86÷6=14, 2"  # was "76÷4=19, 0"
$table.Cell(9, 4).Range.Text = "18÷9=2, 0"  # was "72÷7=10, 2"
$table.Cell(9, 5).Range.Text = "58÷4=14, 2"  # was "50÷9=5, 5"
$table.Cell(13, 1).Range.Text = "44÷3=14, 2"  # was "60÷3=20, 0"
$table.Cell(13, 2).Range.Text = "14÷4=3, 2"  # was "15÷5=3, 0"
$table.Cell(13, 3).Range.Text = "81÷3=27, 0"  # was "43÷5=8, 3"
$table.Cell(13, 4).Range.Text = "79÷5=15, 4"  # was "30÷4=7, 2"
$table.Cell(13, 5).Range.Text = "64÷6=10, 4"  # was "19÷2=9, 1"
$table.Cell(17, 1).Range.Text = "17÷5=3, 2"  # was "12÷6=2, 0"
$table.Cell(17, 2).Range.Text = "92÷3=30, 2"  # was "33÷5=6, 3"
$table.Cell(17, 3).Range.Text = "30÷3=10, 0"  # was "18÷6=3, 0"
$table.Cell(17, 4).Range.Text = "34÷2=17, 0"  # was "70÷8=8, 6"
$table.Cell(17, 5).Range.Text = "72÷4=18, 0"  # was "26÷9=2, 8"
